# "fixed screenshot links compatibility issues with Windows/Linux"
#
# On the "devices" worksheet, the stray leftover test values in K2
# ("Appium") and L2 ("testssts") are removed, and the active selection
# moves to K3 (was L3, now that L2 is empty).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("devices")
$ws.Activate()

$ws.Range("K2:L2").ClearContents()
$ws.Range("K3").Select()
